# Slide 8 ("Potentiometer Connections") originally has, as direct children
# of the slide's shape tree: a picture, the title placeholder, a second
# picture, three arrow connectors, the content placeholder, and a text box.
# The target edit wraps everything except the title and the content
# placeholder into a single new group shape (keeping their relative order
# and geometry untouched), placed where the text box used to be.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# --- Work out what new shape id PowerPoint will hand out -------------------
# A freshly created shape is always assigned the smallest shape id not
# already used on the slide. The shapes being grouped currently carry ids
# up to 36, and native PowerPoint would give the new group id 37 (naming it
# "Group 36"). To reproduce that exactly (instead of whatever low-numbered
# gap happens to be free right now), we burn through every unused id below
# the current maximum first, using disposable textboxes that are deleted
# again immediately. Once those gaps are consumed, the *next* shape created
# is guaranteed to land on maxId + 1.
$maxId = 0
$usedIds = @{}
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $id = $s.Shapes.Item($i).Id
    $usedIds[$id] = $true
    if ($id -gt $maxId) { $maxId = $id }
}

$gapCount = 0
for ($cand = 2; $cand -le $maxId; $cand++) {
    if (-not $usedIds.ContainsKey($cand)) { $gapCount++ }
}

for ($n = 0; $n -lt $gapCount; $n++) {
    $dummy = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
    $dummy.Delete()
}

# --- Build the group ---------------------------------------------------
$namesToGroup = @(
    "Picture 2",
    "Picture 3",
    "Straight Arrow Connector 4",
    "Straight Arrow Connector 5",
    "Straight Arrow Connector 24",
    "TextBox 35"
)

$range = $s.Shapes.Range($namesToGroup)
$grp = $range.Group()

Write-Host "Created $($grp.Name) (id=$($grp.Id)) from $($namesToGroup.Count) shapes"
